$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.959228515625
$ws.Range("B1").Value = 4.852128028869629
$ws.Range("C1").Value = 6.139092445373535
$ws.Range("D1").Value = 10.29194164276123
$ws.Range("E1").Value = 5.179050922393799
